# Refresh the NATMI LR-pair TPM output (Lgi3-Adam22) with newly recomputed
# TPM-based receptor expression / specificity figures (see commit
# "update scripts wuth new tpm"). Only numeric result cells change; the
# dimension columns (A:F, K, L) and labels are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("M2").Value = 8.521337333333333
$ws.Range("N2").Value = 25.564012
$ws.Range("O2").Value = 0.2943426187002489
$ws.Range("P2").Value = 0.2943426187002489
$ws.Range("Q2").Value = 2.552095084200889
$ws.Range("R2").Value = 22.968855757808
$ws.Range("S2").Value = 0.1468608284249291
$ws.Range("T2").Value = 0.1468608284249291
$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("O3").Value = 0.1683364841626613
$ws.Range("P3").Value = 0.1683364841626613
$ws.Range("S3").Value = 0.08399067599328761
$ws.Range("T3").Value = 0.08399067599328763
$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("O4").Value = 0.5373208971370899
$ws.Range("P4").Value = 0.53732089713709
$ws.Range("R4").Value = 41.929524974996
$ws.Range("S4").Value = 0.2680936672780659
$ws.Range("T4").Value = 0.2680936672780661
$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("M5").Value = 8.521337333333333
$ws.Range("N5").Value = 25.564012
$ws.Range("O5").Value = 0.2943426187002489
$ws.Range("P5").Value = 0.2943426187002489
$ws.Range("Q5").Value = 2.562885937710667
$ws.Range("R5").Value = 23.065973439396
$ws.Range("S5").Value = 0.1474817902753197
$ws.Range("T5").Value = 0.1474817902753197
$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("O6").Value = 0.1683364841626613
$ws.Range("P6").Value = 0.1683364841626613
$ws.Range("S6").Value = 0.08434580816937365
$ws.Range("T6").Value = 0.08434580816937366
$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("O7").Value = 0.5373208971370899
$ws.Range("P7").Value = 0.53732089713709
$ws.Range("S7").Value = 0.2692272298590238
$ws.Range("T7").Value = 0.2692272298590239
